$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C28").Value = 326
$ws.Range("D28").Value = 34
$ws.Range("E28").Value = 292
$ws.Range("F28").Value = 5.29595015576324
